# Fix: Use save_sheet instead of undefined save_df
# This restores the rows/sheets that previously failed to persist because
# the (buggy) app code called an undefined save_df instead of save_sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Users: new user "gaston" (row 14)
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

$users.Cells.Item(14, 1).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$users.Cells.Item(14, 2).Value = "gaston"
$users.Cells.Item(14, 3).Value = "gasnmud@gmail.com"
$users.Cells.Item(14, 5).Value = "scrypt:32768:8:1`$a63K5vl6SVSjw0qt`$f44f5ad329399a057bb7f79d0ae944a39ac95424794f6e666da5117acdee2cee40cefa617b9ee7f3c65545e2ef7ce5135a94aa8ef4fed37c269dfed816ed298f"
$users.Cells.Item(14, 6).Value = "ITRADE-28171340"
$users.Cells.Item(14, 7).Value = 0
$users.Cells.Item(14, 9).Value = 45800.44039306713

# Reuse the existing date/time style (s="2") already used by I5:I13, so no
# new style entry is created.
$users.Range("I13").Copy()
$users.Cells.Item(14, 9).PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Beneficiaries: new beneficiary "g" for aaa8d0b9... (row 2)
# ---------------------------------------------------------------------
$beneficiaries = $wb.Worksheets.Item("Beneficiaries")

$beneficiaries.Cells.Item(2, 1).Value = "9d4d8466-2e96-43b3-a2e0-40dcf39823fb"
$beneficiaries.Cells.Item(2, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$beneficiaries.Cells.Item(2, 3).Value = "g"
$beneficiaries.Cells.Item(2, 7).Value = "g"
$beneficiaries.Cells.Item(2, 8).Value = "USD"

# ---------------------------------------------------------------------
# 3) Transactions: four new transactions (rows 3-6) for aaa8d0b9...
# ---------------------------------------------------------------------
$transactions = $wb.Worksheets.Item("Transactions")

# Row 3 - ZAR deposit
$transactions.Cells.Item(3, 1).Value = "9b295e02-a0d7-428f-9533-b95ed325c1a4"
$transactions.Cells.Item(3, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$transactions.Cells.Item(3, 4).Value = 100
$transactions.Cells.Item(3, 7).Value = "ZAR"
$transactions.Cells.Item(3, 8).Value = "Deposit"
$transactions.Cells.Item(3, 9).Value = "Success"
$transactions.Cells.Item(3, 10).Value = 45801.28394535879
$transactions.Range("J2").Copy()
$transactions.Cells.Item(3, 10).PasteSpecial(-4122)

# Row 4 - USD transfer to beneficiary 9d4d8466...
$transactions.Cells.Item(4, 1).Value = "40b27379-5469-4097-a613-a9004267ae8e"
$transactions.Cells.Item(4, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$transactions.Cells.Item(4, 3).Value = "9d4d8466-2e96-43b3-a2e0-40dcf39823fb"
$transactions.Cells.Item(4, 4).Value = 50
$transactions.Cells.Item(4, 6).Value = 45
$transactions.Cells.Item(4, 7).Value = "USD"
$transactions.Cells.Item(4, 8).Value = "Transfer"
$transactions.Cells.Item(4, 9).Value = "Success"
$transactions.Cells.Item(4, 10).Value = 45801.30483710648
$transactions.Range("J2").Copy()
$transactions.Cells.Item(4, 10).PasteSpecial(-4122)

# Row 5 - ZAR deposit
$transactions.Cells.Item(5, 1).Value = "454fdf4d-62ae-46fd-9a75-2b451e94f840"
$transactions.Cells.Item(5, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$transactions.Cells.Item(5, 4).Value = 100
$transactions.Cells.Item(5, 7).Value = "ZAR"
$transactions.Cells.Item(5, 8).Value = "Deposit"
$transactions.Cells.Item(5, 9).Value = "Success"
$transactions.Cells.Item(5, 10).Value = 45801.30553695602
$transactions.Range("J2").Copy()
$transactions.Cells.Item(5, 10).PasteSpecial(-4122)

# Row 6 - USD transfer to beneficiary 9d4d8466...
$transactions.Cells.Item(6, 1).Value = "fcae4e6d-84f2-4bb9-8f1e-48ad94d7bcdb"
$transactions.Cells.Item(6, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$transactions.Cells.Item(6, 3).Value = "9d4d8466-2e96-43b3-a2e0-40dcf39823fb"
$transactions.Cells.Item(6, 4).Value = 100
$transactions.Cells.Item(6, 6).Value = 90
$transactions.Cells.Item(6, 7).Value = "USD"
$transactions.Cells.Item(6, 8).Value = "Transfer"
$transactions.Cells.Item(6, 9).Value = "Success"
$transactions.Cells.Item(6, 10).Value = 45801.30565335074
$transactions.Range("J2").Copy()
$transactions.Cells.Item(6, 10).PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) beneficiaries5: bank_account (G2) was stored as text "62131";
#    fix it up to a real number now that save_sheet persists types
#    correctly.
# ---------------------------------------------------------------------
$beneficiaries5 = $wb.Worksheets.Item("beneficiaries5")
$beneficiaries5.Cells.Item(2, 7).Value = 62131

# ---------------------------------------------------------------------
# 5) New sheet "beneficiaries6" with the first beneficiary for the new
#    "gaston" user (Gerald / bank_account 62312033012 / USD).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$beneficiaries6 = $wb.Worksheets.Add($null, $lastSheet)
$beneficiaries6.Name = "beneficiaries6"

$headerSrc = $beneficiaries5.Range("A1")

$beneficiaries6.Cells.Item(1, 1).Value = "id"
$beneficiaries6.Cells.Item(1, 2).Value = "user_id"
$beneficiaries6.Cells.Item(1, 3).Value = "name"
$beneficiaries6.Cells.Item(1, 4).Value = "id_number"
$beneficiaries6.Cells.Item(1, 5).Value = "country"
$beneficiaries6.Cells.Item(1, 6).Value = "bank_name"
$beneficiaries6.Cells.Item(1, 7).Value = "bank_account"
$beneficiaries6.Cells.Item(1, 8).Value = "currency"
$headerSrc.Copy()
$beneficiaries6.Range("A1:H1").PasteSpecial(-4122)

$beneficiaries6.Cells.Item(2, 1).Value = "e344cb4b-3245-4934-a82a-ed7eee2adf0a"
$beneficiaries6.Cells.Item(2, 2).Value = "aaa8d0b9-9355-4483-ab48-6c0ade3eee45"
$beneficiaries6.Cells.Item(2, 3).Value = "Gerald"
$beneficiaries6.Cells.Item(2, 7).Value = 62312033012
$beneficiaries6.Cells.Item(2, 8).Value = "USD"

# ---------------------------------------------------------------------
# 6) New sheet "LiveRates" with a currency -> USD conversion table.
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$liveRates = $wb.Worksheets.Add($null, $lastSheet2)
$liveRates.Name = "LiveRates"

$liveRates.Cells.Item(1, 1).Value = "currency"
$liveRates.Cells.Item(1, 2).Value = "rate"
$headerSrc.Copy()
$liveRates.Range("A1:B1").PasteSpecial(-4122)

$liveRates.Cells.Item(2, 1).Value = "USD"
$liveRates.Cells.Item(2, 2).Value = 1
$liveRates.Cells.Item(3, 1).Value = "EUR"
$liveRates.Cells.Item(3, 2).Value = 0.881
$liveRates.Cells.Item(4, 1).Value = "GBP"
$liveRates.Cells.Item(4, 2).Value = 0.74
$liveRates.Cells.Item(5, 1).Value = "BWP"
$liveRates.Cells.Item(5, 2).Value = 13.48
$liveRates.Cells.Item(6, 1).Value = "CNY"
$liveRates.Cells.Item(6, 2).Value = 7.18

Write-Output "edit applied"
